$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add Japanese translations in column G (rows 2-22), mirroring column B's English text.
$ws.Range("G2").Value = "弱い声"
$ws.Range("G3").Value = "低鼻音"
$ws.Range("G4").Value = "甲高い泣き声"
$ws.Range("G5").Value = "弱い泣き声"
$ws.Range("G6").Value = "スタッカート泣き"
$ws.Range("G7").Value = "猫鳴き"
$ws.Range("G8").Value = "声の喪失"
$ws.Range("G9").Value = "発声障害"
$ws.Range("G10").Value = "仮性球麻痺"
$ws.Range("G11").Value = "喉頭ジストニア"
$ws.Range("G12").Value = "嗄れ声"
$ws.Range("G13").Value = "異常に低い声"
$ws.Range("G14").Value = "甲高い声"
$ws.Range("G15").Value = "不完全声帯内転"
$ws.Range("G16").Value = "片側声帯麻痺"
$ws.Range("G17").Value = "両側声帯麻痺"
$ws.Range("G18").Value = "両側声帯麻痺"
$ws.Range("G19").Value = "片側声帯麻痺"
$ws.Range("G20").Value = "鼻声"
$ws.Range("G21").Value = "単調発声"
$ws.Range("G22").Value = "音声のピッチ変動の増大"

# Move the active selection to J6, matching the saved view state.
$ws.Range("J6").Select()
